$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row labels: "*_old" -> "*_FV2304", "*_new" -> "*_FV2310" ---
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value()
    if ($v -match '_old$') {
        $cell.Value = ($v -replace '_old$', '_FV2304')
    } elseif ($v -match '_new$') {
        $cell.Value = ($v -replace '_new$', '_FV2310')
    }
}

# --- 2. Turn the data range A1:U58 into an Excel Table (ListObject) named Table1 ---
$tblRange = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add(1, $tblRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row (split/freeze pane below row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
